$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all edited cells keep their original text (string) representation
# instead of being auto-coerced into numbers by the COM Value setter
# (e.g. "1.20" -> 1.2). Force text number format first.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.157.59'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.374.03'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.62%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.99'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.35'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +7.85%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.374.23'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.65%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.40%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.39%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.04%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.21%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.06%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.56%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.384.07'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.13%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.287.83'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.92'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.86'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.66%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '381.78'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +8.03%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.511.23'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.78%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.76'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.16%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +11.39%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +15.75%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +8.17%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.85%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.49%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.10%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.407.79'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.78%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.49%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.49%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.18%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '162.41'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0803'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.87%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.43'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.46'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.71%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.75%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.68'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +8.96%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.20'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.16'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.95'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.09'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +10.90%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.335.98'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.16%  '
